$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.431.89'
$ws.Range('E2').Value = '  -2.60%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.693.60'
$ws.Range('E3').Value = '  -3.22%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '686.24'
$ws.Range('E5').Value = '  -2.34%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '160.97'
$ws.Range('E6').Value = '  -5.94%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.691.87'
$ws.Range('E7').Value = '  -3.25%  '

# Row 8
$ws.Range('E8').Value = '  -0.18%  '

# Row 9
$ws.Range('E9').Value = '  -5.89%  '

# Row 10
$ws.Range('E10').Value = '  -9.20%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.22'
$ws.Range('E11').Value = '  -3.82%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.437'
$ws.Range('E12').Value = '  -10.20%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000234'
$ws.Range('E13').Value = '  -7.03%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.319.54'
$ws.Range('E14').Value = '  -3.29%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '32.65'
$ws.Range('E15').Value = '  -10.67%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.695.06'
$ws.Range('E16').Value = '  -3.42%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '69.519.53'
$ws.Range('E17').Value = '  -2.75%  '

# Row 18
$ws.Range('E18').Value = '  -1.24%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '15.96'
$ws.Range('E19').Value = '  -9.34%  '

# Row 20
$ws.Range('E20').Value = '  -10.67%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '474.15'
$ws.Range('E21').Value = '  -7.86%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.94'
$ws.Range('E22').Value = '  -5.69%  '

# Row 23
$ws.Range('E23').Value = '  -9.52%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.77'
$ws.Range('E24').Value = '  -4.77%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.840.88'
$ws.Range('E25').Value = '  -3.18%  '

# Row 26
$ws.Range('E26').Value = '  -0.02%  '

# Row 27
$ws.Range('E27').Value = '  -11.60%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.05'
$ws.Range('E28').Value = '  -13.57%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.27'
$ws.Range('E29').Value = '  -10.39%  '

# Row 30
$ws.Range('E30').Value = '  -10.00%  '

# Row 31
$ws.Range('E31').Value = '  -11.93%  '

# Row 32
$ws.Range('E32').Value = '  -9.00%  '

# Row 33
$ws.Range('E33').Value = '  -11.17%  '

# Row 34
$ws.Range('E34').Value = '  +0.03%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '26.78'
$ws.Range('E35').Value = '  -8.61%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.160'
$ws.Range('E36').Value = '  -6.47%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '8.23'
$ws.Range('E37').Value = '  -12.19%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.12'
$ws.Range('E38').Value = '  -8.82%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.30'
$ws.Range('E39').Value = '  -3.97%  '

# Row 40
$ws.Range('E40').Value = '  +0.01%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0910'
$ws.Range('E41').Value = '  -10.27%  '

# Row 42
$ws.Range('E42').Value = '  -0.08%  '

# Row 43
$ws.Range('B43').Value = 'Monero'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '167.18'
$ws.Range('E43').Value = '  +0.65%  '

# Row 44
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.944'
$ws.Range('E44').Value = '  -6.77%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '47.88'

# Row 46
$ws.Range('E46').Value = '  -15.34%  '

# Row 47
$ws.Range('E47').Value = '  -3.01%  '

# Row 48
$ws.Range('E48').Value = '  -3.79%  '

# Row 49
$ws.Range('B49').Value = 'FLOKI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.000276'
$ws.Range('E49').Value = '  -9.13%  '

# Row 50
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '28.25'
$ws.Range('E50').Value = '  -7.58%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.85'
$ws.Range('E51').Value = '  -9.29%  '
